# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").Value = "23.020.43"
$ws.Range("E2").Value = "  -3.63%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").Value = "1.601.62"
$ws.Range("E3").Value = "  -2.96%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5: USDC -> USDC
$ws.Range("E5").Value = "  +0.18%  "

# Row 6: BNB -> BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.32"
$ws.Range("E6").Value = "  -3.05%  "

# Row 7: XRP -> XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3783"
$ws.Range("E7").Value = "  -2.83%  "

# Row 8: Cardano -> Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3631"
$ws.Range("E8").Value = "  -5.34%  "

# Row 9: OKB -> OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.09"
$ws.Range("E9").Value = "  -4.24%  "

# Row 10: Polygon -> Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.259"
$ws.Range("E10").Value = "  -6.23%  "

# Row 11: BinanceUSD -> BinanceUSD
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +0.13%  "

# Row 12: Dogecoin -> Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08110"
$ws.Range("E12").Value = "  -3.98%  "

# Row 13: Solana -> Solana
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.77"
$ws.Range("E13").Value = "  -4.59%  "

# Row 14: Polkadot -> Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.586"
$ws.Range("E14").Value = "  -6.20%  "

# Row 15: Chainlink -> Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.409"
$ws.Range("E15").Value = "  -7.27%  "

# Row 16: ShibaInu -> ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001244"
$ws.Range("E16").Value = "  -5.52%  "

# Row 17: WrappedEther -> WrappedEther
$ws.Range("D17").Value = "1.596.01"
$ws.Range("E17").Value = "  -3.26%  "

# Row 18: Litecoin -> Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.94"
$ws.Range("E18").Value = "  -2.22%  "

# Row 19: TRON -> TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06872"
$ws.Range("E19").Value = "  -1.37%  "

# Row 20: Avalanche -> Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.22"
$ws.Range("E20").Value = "  -6.74%  "

# Row 21: Uniswap -> Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.560"
$ws.Range("E21").Value = "  -5.68%  "

# Row 22: Dai -> BitDAO
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5556"
$ws.Range("E22").Value = "  -5.46%  "

# Row 23: Cosmos -> Dai
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.26%  "

# Row 24: WrappedBTC -> Cosmos
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.14"
$ws.Range("E24").Value = "  -3.70%  "

# Row 25: Toncoin -> WrappedBTC
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "23.015.26"
$ws.Range("E25").Value = "  -3.62%  "

# Row 26: LidoDAOToken -> Toncoin
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.378"
$ws.Range("E26").Value = "  -2.63%  "

# Row 27: EthereumClassic -> LidoDAOToken
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.791"
$ws.Range("E27").Value = "  -4.43%  "

# Row 28: Monero -> EthereumClassic
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.07"
$ws.Range("E28").Value = "  -4.12%  "

# Row 29: HuobiToken -> Monero
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.47"
$ws.Range("E29").Value = "  -2.35%  "

# Row 30: BitcoinCash -> HuobiToken
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.246"
$ws.Range("E30").Value = "  -2.49%  "

# Row 31: WEMIXTOKEN -> BitcoinCash
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.23"
$ws.Range("E31").Value = "  -2.99%  "

# Row 32: Filecoin -> WEMIXTOKEN
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.307"
$ws.Range("E32").Value = "  -7.24%  "

# Row 33: WrappedliquidstakedEther2.0 -> Filecoin
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.839"
$ws.Range("E33").Value = "  -11.71%  "

# Row 34: ImmutableX -> WrappedliquidstakedEther2.0
$ws.Range("B34").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C34").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D34").Value = "1.782.62"
$ws.Range("E34").Value = "  -2.70%  "

# Row 35: Hedera -> ImmutableX
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9620"
$ws.Range("E35").Value = "  -3.00%  "

# Row 36: FraxShare -> Hedera
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07632"
$ws.Range("E36").Value = "  -6.46%  "

# Row 37: InternetComputer(DFINITY) -> FraxShare
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.36"
$ws.Range("E37").Value = "  -1.64%  "

# Row 38: VeChain -> InternetComputer(DFINITY)
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.281"
$ws.Range("E38").Value = "  -5.89%  "

# Row 39: Algorand -> VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02711"
$ws.Range("E39").Value = "  -6.76%  "

# Row 40: Stellar -> Algorand
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2530"
$ws.Range("E40").Value = "  -5.44%  "

# Row 41: TrustWalletToken -> Stellar
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08843"
$ws.Range("E41").Value = "  -2.87%  "

# Row 42: TheSandbox -> TrustWalletToken
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.364"
$ws.Range("E42").Value = "  -4.03%  "

# Row 43: Aptos -> TheSandbox
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7044"
$ws.Range("E43").Value = "  -6.75%  "

# Row 44: EnergySwap -> Aptos
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.53"
$ws.Range("E44").Value = "  -7.21%  "

# Row 45: Decentraland -> EnergySwap
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.21"
$ws.Range("E45").Value = "  -9.01%  "

# Row 46: Frax -> Decentraland
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6604"
$ws.Range("E46").Value = "  -4.74%  "

# Row 47: NEARProtocol -> Frax
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  +0.15%  "

# Row 48: PancakeSwap -> NEARProtocol
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.315"
$ws.Range("E48").Value = "  -5.46%  "

# Row 49: Quant -> PancakeSwap
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.995"
$ws.Range("E49").Value = "  -2.46%  "

# Row 50: Cronos -> Quant
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.48"
$ws.Range("E50").Value = "  -1.15%  "

# Row 51: Flow -> Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07914"
$ws.Range("E51").Value = "  -4.26%  "

